$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B4:B73 values per the naive forecaster bugfix
$ws.Range("B4").Value = -0.3499999999999943
$ws.Range("B5").Value = 0.25
$ws.Range("B6").Value = -0.8999999999999915
$ws.Range("B7").Value = -0.6500000000000057
$ws.Range("B8").Value = -0.7000000000000171
$ws.Range("B9").Value = 0.2000000000000028
$ws.Range("B10").Value = 0.4999999999999858
$ws.Range("B11").Value = 0.4200000000000017
$ws.Range("B12").Value = 1.099999999999994
$ws.Range("B13").Value = 0.4999999999999858
$ws.Range("B14").Value = 0.8
$ws.Range("B15").Value = 0.3
$ws.Range("B16").Value = 0.4
$ws.Range("B17").Value = 0.5
$ws.Range("B18").Value = -0.2
$ws.Range("B19").Value = 0.1
$ws.Range("B20").Value = 0.1
$ws.Range("B21").Value = 0.1
$ws.Range("B22").Value = -0.3
$ws.Range("B23").Value = 0.2
$ws.Range("B24").Value = 0.8999999999999915
$ws.Range("B25").Value = 0.4200000000000017
$ws.Range("B26").Value = 0.3
$ws.Range("B27").Value = 0.5
$ws.Range("B28").Value = 0.2999999999999829
$ws.Range("B29").Value = 0.4000000000000057
$ws.Range("B30").Value = 0.2000000000000028
$ws.Range("B31").Value = 0.4200000000000017
$ws.Range("B32").Value = 0.6200000000000045
$ws.Range("B33").Value = 0.4200000000000017
$ws.Range("B34").Value = 0.3000000000000114
$ws.Range("B35").Value = 0.5400000000000063
$ws.Range("B36").Value = 0.3400000000000034
$ws.Range("B37").Value = 0.4399999999999977
$ws.Range("B38").Value = 0.4999999999999858
$ws.Range("B39").Value = 0.4999999999999858
$ws.Range("B40").Value = 0.6999999999999886
$ws.Range("B41").Value = 0.5999999999999943
$ws.Range("B42").Value = 0.5999999999999943
$ws.Range("B43").Value = 0.6999999999999886
$ws.Range("B44").Value = 0.2999999999999829
$ws.Range("B45").Value = 0.4999999999999716
$ws.Range("B46").Value = 0.2
$ws.Range("B47").Value = 0.09999999999999432
$ws.Range("B48").Value = -0.1
$ws.Range("B49").Value = -0.09999999999999432
$ws.Range("B50").Value = 0.08000000000004093
$ws.Range("B51").Value = 0.09999999999999432
$ws.Range("B52").Value = -11.9
$ws.Range("B53").Value = 6.640000000000001
$ws.Range("B54").Value = -0.4000000000000057
$ws.Range("B55").Value = -0.7094799999999992
$ws.Range("B56").Value = 1.310000000000016
$ws.Range("B57").Value = 1.52000000000001
$ws.Range("B58").Value = -0.539999999999992
$ws.Range("B59").Value = 0.4652855479103435
$ws.Range("B60").Value = 0.38
$ws.Range("B61").Value = 0.04
$ws.Range("B62").Value = -0.29
$ws.Range("B63").Value = -0.2078779574152918
$ws.Range("B64").Value = 0.1206478331785803
$ws.Range("B65").Value = -0.18
$ws.Range("B66").Value = 0.044
$ws.Range("B67").Value = -0.08251004046350374
$ws.Range("B68").Value = 0.2582525219575302
$ws.Range("B69").Value = -0.04717552522494373
$ws.Range("B70").Value = 0.2142297805489477
$ws.Range("B71").Value = 0.2394371574146135
$ws.Range("B72").Value = 0.04717883418304325
$ws.Range("B73").Value = 0.0959495356205764

# Remove rows 74-82 (data series shortened)
$ws.Range("A74:B82").EntireRow.Delete() | Out-Null

$ws.Range("A1:B73").Select() | Out-Null
